$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 0
    9  = 2
    10 = 2
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
